# Brockhaus AG link von Folienmaster entfernt
#
# The "Rectangle 39" shape that carries the hyperlink to
# http://www.brockhaus-ag.de/ is a userDrawn shape living directly on the
# slide masters (not on individual slides/layouts). It is present on two of
# the presentation's three slide masters (the ones that back slides 2-3 and
# slides 4-15 respectively) - remove it from every slide master that has it.

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Designs.Count; $i++) {
    $master = $p.Designs.Item($i).SlideMaster

    for ($j = $master.Shapes.Count; $j -ge 1; $j--) {
        $shp = $master.Shapes.Item($j)
        if ($shp.Name -eq "Rectangle 39") {
            Write-Host "Design" $i ": removing shape '" $shp.Name "' (id" $shp.Id ") from slide master"
            $shp.Delete()
        }
    }
}
